$wb = $excel.ActiveWorkbook

# --- Sheet: Cost of Production ---
$ws3 = $wb.Worksheets.Item("Cost of Production")

$ws3.Range("E2").Value = 92585.2570508088
$ws3.Range("F2").Value = 249632.057050809
$ws3.Range("G2").Value = 249632.057050809

$ws3.Range("E3").Value = 93084.9820508088
$ws3.Range("F3").Value = 109225.382050809
$ws3.Range("G3").Value = 358857.439101618

$ws3.Range("E4").Value = 103027.008300809
$ws3.Range("F4").Value = 208873.289550809
$ws3.Range("G4").Value = 567730.728652426

$ws3.Range("D5").Value = 35157.78125
$ws3.Range("E5").Value = 104069.034550809
$ws3.Range("F5").Value = 139226.815800809
$ws3.Range("G5").Value = 706957.544453235
$ws3.Range("H5").Value = 0
$ws3.Range("I5").Value = 0
$ws3.Range("J5").Value = 0
$ws3.Range("K5").Value = 0
$ws3.Range("L5").Value = 0
$ws3.Range("M5").Value = 0
$ws3.Range("N5").Value = "#NUM!"
$ws3.Range("O5").Value = "#NUM!"
$ws3.Range("P5").Value = "#NUM!"
$ws3.Range("Q5").Value = "#NUM!"

$ws3.Range("D6").Value = 41675
$ws3.Range("F6").Value = 145744.034550809
$ws3.Range("G6").Value = 852701.579004044
$ws3.Range("H6").Value = 58618.1640625
$ws3.Range("I6").Value = 120.321920402265
$ws3.Range("J6").Value = 4.73708603842923
$ws3.Range("K6").Value = 42.5074125817755
$ws3.Range("L6").Value = 0.0937126919260338
$ws3.Range("M6").Value = 10.6709131863301
$ws3.Range("N6").Value = 2.48632888596465
$ws3.Range("O6").Value = 14.5467124848004
$ws3.Range("P6").Value = 26.5313996947935
$ws3.Range("Q6").Value = 155.226706071808

$ws3.Range("D7").Value = 41675
$ws3.Range("F7").Value = 145744.034550809
$ws3.Range("G7").Value = 998445.613554853
$ws3.Range("H7").Value = 58618.1640625
$ws3.Range("I7").Value = 120.321920402265
$ws3.Range("J7").Value = 4.73708603842923
$ws3.Range("K7").Value = 42.5074125817755
$ws3.Range("L7").Value = 0.0937126919260338
$ws3.Range("M7").Value = 10.6709131863301
$ws3.Range("N7").Value = 2.48632888596465
$ws3.Range("O7").Value = 8.5165206853825
$ws3.Range("P7").Value = 26.5313996947935
$ws3.Range("Q7").Value = 90.8790528833008

$ws3.Range("D8").Value = 41675
$ws3.Range("F8").Value = 145744.034550809
$ws3.Range("G8").Value = 1144189.64810566
$ws3.Range("H8").Value = 58618.1640625
$ws3.Range("I8").Value = 120.321920402265
$ws3.Range("J8").Value = 4.73708603842923
$ws3.Range("K8").Value = 42.5074125817755
$ws3.Range("L8").Value = 0.0937126919260338
$ws3.Range("M8").Value = 10.6709131863301
$ws3.Range("N8").Value = 2.48632888596465
$ws3.Range("O8").Value = 6.50645675224322
$ws3.Range("P8").Value = 26.5313996947935
$ws3.Range("Q8").Value = 69.4298351537983

$ws3.Range("D9").Value = 41675
$ws3.Range("F9").Value = 145744.034550809
$ws3.Range("G9").Value = 1289933.68265647
$ws3.Range("H9").Value = 58618.1640625
$ws3.Range("I9").Value = 120.321920402265
$ws3.Range("J9").Value = 4.73708603842923
$ws3.Range("K9").Value = 42.5074125817755
$ws3.Range("L9").Value = 0.0937126919260338
$ws3.Range("M9").Value = 10.6709131863301
$ws3.Range("N9").Value = 2.48632888596465
$ws3.Range("O9").Value = 5.50142478567358
$ws3.Range("P9").Value = 26.5313996947935
$ws3.Range("Q9").Value = 58.7052262890471

$ws3.Range("D10").Value = 41675
$ws3.Range("F10").Value = 145744.034550809
$ws3.Range("G10").Value = 1435677.71720728
$ws3.Range("H10").Value = 58618.1640625
$ws3.Range("I10").Value = 120.321920402265
$ws3.Range("J10").Value = 4.73708603842923
$ws3.Range("K10").Value = 42.5074125817755
$ws3.Range("L10").Value = 0.0937126919260338
$ws3.Range("M10").Value = 10.6709131863301
$ws3.Range("N10").Value = 2.48632888596465
$ws3.Range("O10").Value = 4.89840560573179
$ws3.Range("P10").Value = 26.5313996947935
$ws3.Range("Q10").Value = 52.2704609701964

$ws3.Range("D11").Value = 41675
$ws3.Range("F11").Value = 145744.034550809
$ws3.Range("G11").Value = 1581421.75175809
$ws3.Range("H11").Value = 58618.1640625
$ws3.Range("I11").Value = 120.321920402265
$ws3.Range("J11").Value = 4.73708603842923
$ws3.Range("K11").Value = 42.5074125817755
$ws3.Range("L11").Value = 0.0937126919260338
$ws3.Range("M11").Value = 10.6709131863301
$ws3.Range("N11").Value = 2.48632888596465
$ws3.Range("O11").Value = 4.49639281910393
$ws3.Range("P11").Value = 26.5313996947935
$ws3.Range("Q11").Value = 47.9806174242959

$ws3.Range("D12").Value = 41675
$ws3.Range("F12").Value = 145744.034550809
$ws3.Range("G12").Value = 1727165.7863089
$ws3.Range("H12").Value = 58618.1640625
$ws3.Range("I12").Value = 120.321920402265
$ws3.Range("J12").Value = 4.73708603842923
$ws3.Range("K12").Value = 42.5074125817755
$ws3.Range("L12").Value = 0.0937126919260338
$ws3.Range("M12").Value = 10.6709131863301
$ws3.Range("N12").Value = 2.48632888596465
$ws3.Range("O12").Value = 4.20924082865547
$ws3.Range("P12").Value = 26.5313996947935
$ws3.Range("Q12").Value = 44.9164434629384


# --- Sheet: Labor (insert new row for Dropper Line Cleaning (Spring), Y3; shift subsequent rows) ---
$ws5 = $wb.Worksheets.Item("Labor")

$ws5.Range("A16").Value = "15"
$ws5.Range("B16").Value = "Dropper Line Cleaning (Spring)"
$ws5.Range("C16").Value = "Y3"
$ws5.Range("D16").Value = "Ear Hanging"
$ws5.Range("E16").Value = "Y"
$ws5.Range("F16").Value = 15.9127272727273
$ws5.Range("G16").Value = 2
$ws5.Range("H16").Value = "Devices/Day"
$ws5.Range("I16").Value = 275
$ws5.Range("J16").Value = 1
$ws5.Range("K16").Value = "Spring"
$ws5.Range("L16").Value = "Dropper line cleaning, important for growth and `nshell quality in whole scallop market.  Cleaning `nrequires a specialized system to do easily"
$ws5.Range("M16").Value = 16
$ws5.Range("N16").Value = 240

$ws5.Range("A17").Value = "16"
$ws5.Range("B17").Value = "Dropper Line Cleaning (Summer)"
$ws5.Range("C17").Value = "Y3"
$ws5.Range("D17").Value = "Ear Hanging"
$ws5.Range("E17").Value = "Y"
$ws5.Range("F17").Value = 15.9127272727273
$ws5.Range("G17").Value = 2
$ws5.Range("H17").Value = "Devices/Day"
$ws5.Range("I17").Value = 275
$ws5.Range("J17").Value = 1
$ws5.Range("K17").Value = "Summer"
$ws5.Range("L17").Value = "Dropper line cleaning, important for growth and `nshell quality in whole scallop market.  Cleaning `nrequires a specialized system to do easily"
$ws5.Range("M17").Value = 16
$ws5.Range("N17").Value = 240

$ws5.Range("A18").Value = "17"
$ws5.Range("B18").Value = "Dropper Line Cleaning (Fall)"
$ws5.Range("C18").Value = "Y3"
$ws5.Range("D18").Value = "Ear Hanging"
$ws5.Range("E18").Value = "Y"
$ws5.Range("F18").Value = 15.9127272727273
$ws5.Range("G18").Value = 2
$ws5.Range("H18").Value = "Devices/Day"
$ws5.Range("I18").Value = 275
$ws5.Range("J18").Value = 1
$ws5.Range("K18").Value = "Fall"
$ws5.Range("L18").Value = "Dropper line cleaning, important for growth and `nshell quality in whole scallop market.  Cleaning `nrequires a specialized system to do easily"
$ws5.Range("M18").Value = 16
$ws5.Range("N18").Value = 240

$ws5.Range("A19").Value = "18"
$ws5.Range("B19").Value = "Harvest"
$ws5.Range("C19").Value = "Y3"
$ws5.Range("D19").Value = "Global"
$ws5.Range("E19").Value = "Y"
$ws5.Range("F19").Value = 375.15625
$ws5.Range("G19").Value = 47
$ws5.Range("H19").Value = "Scallops/Day"
$ws5.Range("I19").Value = 1250
$ws5.Range("J19").Value = 1
$ws5.Range("K19").Value = "Fall"
$ws5.Range("L19").Value = "Harvest, similar time frame for whole scallops`nand adductor only.  Does not account for`ndelivery to market"
$ws5.Range("M19").Value = 376
$ws5.Range("N19").Value = 5640

$ws5.Range("A20").Value = "19"
$ws5.Range("B20").Value = "Prep Time (Summer)"
$ws5.Range("C20").Value = "all"
$ws5.Range("D20").Value = "Global"
$ws5.Range("E20").Value = "Y"
$ws5.Range("F20").Value = 120
$ws5.Range("G20").Value = 0
$ws5.Range("H20").Value = "Hours/Week"
$ws5.Range("I20").Value = 10
$ws5.Range("J20").Value = 0
$ws5.Range("K20").Value = "Summer"
$ws5.Range("L20").Value = "On land 'Prep time' for misc tasks `nex: cement bucket anchor construction"
$ws5.Range("M20").Value = 120
$ws5.Range("N20").Value = 0

$ws5.Range("A21").Value = "20"
$ws5.Range("B21").Value = "Prep Time (Fall)"
$ws5.Range("C21").Value = "all"
$ws5.Range("D21").Value = "Global"
$ws5.Range("E21").Value = "Y"
$ws5.Range("F21").Value = 0
$ws5.Range("G21").Value = 0
$ws5.Range("H21").Value = "Hours/Week"
$ws5.Range("I21").Value = 0
$ws5.Range("J21").Value = 0
$ws5.Range("K21").Value = "Fall"
$ws5.Range("L21").Value = "On land 'Prep time' for misc tasks `nex: cement bucket anchor construction"
$ws5.Range("M21").Value = 0
$ws5.Range("N21").Value = 0

$ws5.Range("A22").Value = "21"
$ws5.Range("B22").Value = "Prep Time (Spring)"
$ws5.Range("C22").Value = "all"
$ws5.Range("D22").Value = "Global"
$ws5.Range("E22").Value = "Y"
$ws5.Range("F22").Value = 0
$ws5.Range("G22").Value = 0
$ws5.Range("H22").Value = "Hours/Week"
$ws5.Range("I22").Value = 0
$ws5.Range("J22").Value = 0
$ws5.Range("K22").Value = "Spring"
$ws5.Range("L22").Value = "On land 'Prep time' for misc tasks `nex: cement bucket anchor construction"
$ws5.Range("M22").Value = 0
$ws5.Range("N22").Value = 0

$ws5.Range("A23").Value = "22"
$ws5.Range("B23").Value = "Prep Time (Winter)"
$ws5.Range("C23").Value = "all"
$ws5.Range("D23").Value = "Global"
$ws5.Range("E23").Value = "Y"
$ws5.Range("F23").Value = 0
$ws5.Range("G23").Value = 0
$ws5.Range("H23").Value = "Hours/Week"
$ws5.Range("I23").Value = 0
$ws5.Range("J23").Value = 0
$ws5.Range("K23").Value = "Winter"
$ws5.Range("L23").Value = "On land 'Prep time' for misc tasks `nex: cement bucket anchor construction"
$ws5.Range("M23").Value = 0
$ws5.Range("N23").Value = 0


# --- Sheet: Fuel ---
$ws6 = $wb.Worksheets.Item("Fuel")
$ws6.Range("I8").Value = 6678
$ws6.Range("I9").Value = 3180

# --- Sheet: Maintenance ---
$ws7 = $wb.Worksheets.Item("Maintenance")
$ws7.Range("E9").Value = 7950
$ws7.Range("E10").Value = 79.5
